$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell explicitly so text-like numeric strings
# (prices / percentages) are preserved verbatim as text, not coerced
# into floating point numbers by the normal Value setter. Forcing the
# NumberFormat to "@" (Text) before assigning keeps the literal string,
# then resetting the style back to Normal removes the temporary format.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '297.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-1.23%'
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '31.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '0.67%'
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-1.52%'
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08020'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '8.71%'
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.655'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '28.08%'
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '7.788'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-1.32%'
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.821'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '1.59%'
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9257'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.31%'
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1755'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.72%'
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07406'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.42%'
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08956'
$ws.Range('D12').Style = 'Normal'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03066'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.83%'
$ws.Range('E13').Style = 'Normal'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001490'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.31%'
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005872'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-3.92%'
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.538'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '2.62%'
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.248'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.89%'
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.11%'
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1336'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.04%'
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.044'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-12.98%'
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '4.04%'
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04605'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-1.08%'
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001243'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.93%'
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004429'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.99%'
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-7.69%'
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003411'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '82.24%'
$ws.Range('E27').Style = 'Normal'

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01761'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '2.33%'
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04498'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.33%'
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006700'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-6.23%'
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1347'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.15%'
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002209'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '3.76%'
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009845'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-6.11%'
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00006487'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '3.48%'
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.06%'
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('B47').Value = 'BOLO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.8206'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '11.20%'
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.008745'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '24.99%'
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002099'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.06%'
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.01%'
$ws.Range('E50').Style = 'Normal'
